$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 3 now carries what used to be the "AC_course_support_staff" test case
# (previously on row 9), and row 4 now carries what used to be the
# "AC_teaching staff" test case (previously on row 10).
$ws1.Range("A3").Value = "AC_course_support_staff"
$ws1.Range("B3").Value = "aa000fyl"
$ws1.Range("C3").Value = "Course Support Staff"
$ws1.Range("D3").Value = 11

$ws1.Range("A4").Value = "AC_teaching staff "
$ws1.Range("B4").Value = "aa000fzm"
$ws1.Range("C4").Value = "Teaching Staff"
$ws1.Range("D4").Value = 12

# The rest of the old test rows (5-12) are no longer needed - clear their
# contents but keep the existing row/cell styling in place.
$ws1.Range("A5:D12").ClearContents()

# Leave the active selection on A11, matching where the user ended up after
# trimming the sheet down.
[void]$ws1.Activate()
[void]$ws1.Range("A11").Select()
